# Fix mislabeled columns in the SectorGroup codelist sheet.
#
# The sheet has 7 columns: code, name, status, category-name, group-name,
# group-code, category-code. The category/group name columns (D/E) and the
# group/category code columns (F/G) were swapped by mistake, so this script
# swaps column D with column E, and column F with column G, for every row
# (including the header row), to restore the correct mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Columns F and G hold numeric-looking codes (e.g. "110") that must stay
# text, so force a text number format before writing back into them.
$ws.Range("F1:G" + $lastRow).NumberFormat = "@"

for ($r = 1; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2
    $gVal = $ws.Cells.Item($r, 7).Value2

    $ws.Cells.Item($r, 4).Value2 = $eVal
    $ws.Cells.Item($r, 5).Value2 = $dVal
    $ws.Cells.Item($r, 6).Value2 = $gVal
    $ws.Cells.Item($r, 7).Value2 = $fVal
}
